# Update the per-region breakdown figures ("davon": Auslaenderinnen und
# Auslaender / (Spaet-)Aussiedlerinnen und (Spaet-)Aussiedler / weitere
# Deutsche mit Migrationshintergrund) for 2018 in columns K:M, rows 27-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 27;  K = 103.07277000000001;   L = 44.397529999999996;   M = 93.498339999999999 },
    @{ Row = 28;  K = 42.428899999999999;   L = 19.416730000000001;   M = 42.408319999999996 },
    @{ Row = 29;  K = 145.50167000000002;   L = 63.814260000000004;   M = 135.90667000000002 },
    @{ Row = 30;  K = 159.72075000000001;   L = 39.890329999999999;   M = 147.85239999999999 },
    @{ Row = 31;  K = 91.846050000000005;   L = 20.654400000000003;   M = 89.197130000000001 },
    @{ Row = 32;  K = 67.87469999999999;    L = 19.23593;             M = 58.655269999999994 },
    @{ Row = 33;  K = 43.481569999999998;   L = 22.718240000000002;   M = 43.496079999999999 },
    @{ Row = 34;  K = 38.967309999999998;   L = 14.056089999999999;   M = 35.339649999999999 },
    @{ Row = 35;  K = 242.16964000000002;   L = 76.664649999999995;   M = 226.68813 },
    @{ Row = 36;  K = 65.021280000000004;   L = 19.552619999999997;   M = 77.539289999999994 },
    @{ Row = 37;  K = 56.90699;             L = 18.322279999999999;   M = 57.331009999999999 },
    @{ Row = 38;  K = 121.92827;            L = 37.874900000000004;   M = 134.87029999999999 },
    @{ Row = 39;  K = 54.569319999999998;   L = 15.27867;             M = 44.010460000000002 },
    @{ Row = 40;  K = 51.998440000000002;   L = 30.461880000000001;   M = 58.261160000000004 },
    @{ Row = 41;  K = 134.04888;            L = 58.97766;             M = 93.08153999999999 },
    @{ Row = 42;  K = 240.61664000000002;   L = 104.71821000000001;   M = 195.35316 },
    @{ Row = 43;  K = 750.21622000000002;   L = 283.07202000000001;   M = 692.81825000000003 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 11).Value = $u.K   # column K
    $ws.Cells.Item($r, 12).Value = $u.L   # column L
    $ws.Cells.Item($r, 13).Value = $u.M   # column M
}

# Restore the view/cursor state: scroll so row 19 is the top visible row and
# leave the active cell on J35, matching where the author was working.
[void]$excel.Goto($ws.Range("A19"), $true)
[void]$ws.Range("J35").Select()
